# Fix the extra backslash in the developer/tester names on the
# "ProductsData" sheet (column F: Bharadwaj, column G: jayaraj), and
# restore the selection to the range that was being edited (G2:G11).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductsData")

for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 6).Value = "valgenesis\Bharadwaj"
    $ws.Cells.Item($r, 7).Value = "valgenesis\jayaraj"
}

$ws.Activate()
$ws.Range("G2:G11").Select()
